$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = -0.1112181501777582
$ws.Range("C2").Value = 2.129230451805604
$ws.Range("D2").Value = 21.3896688524574
$ws.Range("E2").Value = 4.624896631542958
$ws.Range("F2").Value = 4.732363681339207
$ws.Range("G2").Value = 22

# Row 3 (Q1)
$ws.Range("B3").Value = 1.601564424175498
$ws.Range("C3").Value = 2.22790675060973
$ws.Range("D3").Value = 23.21886171173811
$ws.Range("E3").Value = 4.818595408595549
$ws.Range("F3").Value = 4.656881549094973
$ws.Range("G3").Value = 21

# Row 4 (Q2)
$ws.Range("B4").Value = 0.2290663667550482
$ws.Range("C4").Value = 1.376689112904139
$ws.Range("D4").Value = 4.866116910555858
$ws.Range("E4").Value = 2.205927675730974
$ws.Range("F4").Value = 2.250998728091834
$ws.Range("G4").Value = 20

# Row 5 (Q3)
$ws.Range("B5").Value = 0.7581909808868581
$ws.Range("C5").Value = 0.9427861711352807
$ws.Range("D5").Value = 1.287187669203633
$ws.Range("E5").Value = 1.134542934050375
$ws.Range("F5").Value = 0.8671264167865568
$ws.Range("G5").Value = 19

# Row 6 (Q4)
$ws.Range("B6").Value = 0.6536748506019255
$ws.Range("C6").Value = 0.8947687390269206
$ws.Range("D6").Value = 1.324266432758395
$ws.Range("E6").Value = 1.150767757959179
$ws.Range("F6").Value = 0.9745454809077441
$ws.Range("G6").Value = 18

# Row 7 (Q5)
$ws.Range("B7").Value = 0.2081861679188835
$ws.Range("C7").Value = 0.6178496757196613
$ws.Range("D7").Value = 0.5151342045300913
$ws.Range("E7").Value = 0.7177285033563118
$ws.Range("F7").Value = 0.7080111364014168
$ws.Range("G7").Value = 17

# Row 8 (Q6)
$ws.Range("B8").Value = 0.1864016188994063
$ws.Range("C8").Value = 0.5392751603362338
$ws.Range("D8").Value = 0.464854709378989
$ws.Range("E8").Value = 0.6818025442743588
$ws.Range("F8").Value = 0.6773352854439084
$ws.Range("G8").Value = 16
